$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (laptop) "facts" cell B2: updated product description text
# (rewording / reflow of the sentence around the "\n" marker before the
# carbon-footprint line).
$newFacts = 'Ekran o przekątnej 15.6" i wysokiej rozdzielczości z matową powłoką ograniczającą odbijanie się promieni słonecznych\n 
Dwurdzeniowy procesor i 8 GB pamięci RAM  pozwalającej na uruchomienie kilku aplikacji jednocześnie oraz granie w wymagające sprzętowo gry komputerowe, pojemność dysku SSD 256 GB.\n
Ślad węglowy to 423 kg'
$ws.Range("B2").Value = $newFacts

# Row 5 (smartband) "image" cell E5: point to the new figure file
$ws.Range("E5").Value = "png/smartband2.png"

# Leave the selection on the cell that was just edited (matches the
# author's final cursor position after updating the image reference).
$ws.Range("E5").Select() | Out-Null
